# Update the "dSF" column (F) values for a set of rows to reflect
# re-pulled data. Row numbers below are the Excel worksheet row numbers
# (1-based), matching the XML <row r="N"> of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    4  = -2
    9  = 7
    14 = 3
    18 = -6
    21 = 1
    24 = -8
    26 = -4
    27 = -3
    28 = -5
    34 = -4
    36 = -1
    37 = 0
    39 = -2
    40 = -2
    41 = 3
    47 = 3
    50 = -3
    52 = -1
    55 = 1
    60 = -2
    72 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
